$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text so numeric-looking values (e.g. "0.9997")
# are stored as text, matching the source data feed formatting.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.230.48'
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").Value = '1.863.15'

$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '0.7093'
$ws.Range("E5").Value = '  +0.89%  '

$ws.Range("D6").Value = '237.75'
$ws.Range("E6").Value = '  -0.32%  '

$ws.Range("D7").Value = '0.9994'
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").Value = '0.08206'
$ws.Range("E8").Value = '  +9.56%  '

$ws.Range("E9").Value = '  -0.69%  '

$ws.Range("D10").Value = '23.32'
$ws.Range("E10").Value = '  -0.58%  '

$ws.Range("D11").Value = '0.08171'
$ws.Range("E11").Value = '  +0.57%  '

$ws.Range("D12").Value = '1.820.28'
$ws.Range("E12").Value = '  -2.84%  '

$ws.Range("D13").Value = '5.168'
$ws.Range("E13").Value = '  -1.07%  '

$ws.Range("D14").Value = '0.7086'
$ws.Range("E14").Value = '  -2.42%  '

$ws.Range("D15").Value = '89.30'
$ws.Range("E15").Value = '  +0.55%  '

$ws.Range("D16").Value = '29.225.83'
$ws.Range("E16").Value = '  -0.30%  '

$ws.Range("D17").Value = '0.000007898'
$ws.Range("E17").Value = '  +3.54%  '

$ws.Range("D18").Value = '5.786'
$ws.Range("E18").Value = '  +0.31%  '

$ws.Range("D19").Value = '13.34'
$ws.Range("E19").Value = '  +2.17%  '

$ws.Range("D20").Value = '237.12'
$ws.Range("E20").Value = '  -0.67%  '

$ws.Range("E21").Value = '  -0.10%  '

$ws.Range("D22").Value = '2.104.88'
$ws.Range("E22").Value = '  -1.01%  '

$ws.Range("D23").Value = '0.9995'
$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").Value = '7.404'
$ws.Range("E24").Value = '  -2.38%  '

$ws.Range("D25").Value = '162.21'
$ws.Range("E25").Value = '  +0.67%  '

$ws.Range("D26").Value = '8.965'
$ws.Range("E26").Value = '  -0.38%  '

$ws.Range("D27").Value = '0.1444'
$ws.Range("E27").Value = '  -0.70%  '

$ws.Range("D28").Value = '18.08'
$ws.Range("E28").Value = '  -0.01%  '

$ws.Range("E29").Value = '  -1.40%  '

$ws.Range("E30").Value = '  +1.98%  '

$ws.Range("D31").Value = '1.482'
$ws.Range("E31").Value = '  -0.75%  '

$ws.Range("D32").Value = '4.387'
$ws.Range("E32").Value = '  -3.55%  '

$ws.Range("D33").Value = '4.053'
$ws.Range("E33").Value = '  +1.94%  '

$ws.Range("D34").Value = '0.05208'
$ws.Range("E34").Value = '  +0.71%  '

$ws.Range("D35").Value = '1.170'
$ws.Range("E35").Value = '  -1.41%  '

$ws.Range("D36").Value = '0.7071'
$ws.Range("E36").Value = '  +0.57%  '

$ws.Range("D37").Value = '0.9990'
$ws.Range("E37").Value = '  -3.74%  '

$ws.Range("D38").Value = '2.669'
$ws.Range("E38").Value = '  +0.58%  '

$ws.Range("D39").Value = '0.01847'
$ws.Range("E39").Value = '  -1.02%  '

$ws.Range("D40").Value = '2.727'
$ws.Range("E40").Value = '  +1.87%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '0.9246'
$ws.Range("E41").Value = '  -0.71%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.142.71'
$ws.Range("E42").Value = '  +6.28%  '

$ws.Range("D43").Value = '0.4275'
$ws.Range("E43").Value = '  -0.18%  '

$ws.Range("D44").Value = '5.866'
$ws.Range("E44").Value = '  -2.53%  '

$ws.Range("D45").Value = '70.16'
$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("D46").Value = '0.9987'
$ws.Range("E46").Value = '  -0.07%  '

$ws.Range("D47").Value = '102.71'
$ws.Range("E47").Value = '  +0.03%  '

$ws.Range("D48").Value = '1.776'
$ws.Range("E48").Value = '  +1.82%  '

$ws.Range("E49").Value = '  -0.68%  '

$ws.Range("D50").Value = '9.196'
$ws.Range("E50").Value = '  +0.47%  '

$ws.Range("D51").Value = '6.958'
$ws.Range("E51").Value = '  -1.12%  '

# Restore default (General) formatting now that the text values are committed,
# so no residual text-number-format style is left on the Price column.
$priceRange.ClearFormats()